# Add the "3b" variant result rows (periodMin = 8, rather than no minimum
# period) to the summary results table: "C naive3b" (C version) and
# "P naive3b" (P version). This mirrors re-running the benchmark suite with
# the amended 3b algorithm and pasting its row into the results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6; this pushes the existing rows 6-9
# ("P algorithm", "P naive1", "P naive2", "P naive3") down to rows 7-10,
# leaving row 6 free for the new "C naive3b" result row right after the
# other "C ..." rows (rows 2-5).
$ws.Rows("6:6").Insert()

# Row 6: C naive3b
$ws.Range("A6").Value = "C naive3b"
$ws.Range("B6").Value = 316
$ws.Range("C6").Value = 48
$ws.Range("D6").Value = 95
$ws.Range("E6").Value = 95
$ws.Range("F6").Value = 97
$ws.Range("G6").Value = 428
$ws.Range("H6").Value = 137.0601265822785
$ws.Range("I6").Value = 120.1679757883243
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 4
$ws.Range("L6").Value = 5
$ws.Range("M6").Value = 6
$ws.Range("N6").Value = 38
$ws.Range("O6").Value = 7.981012658227848
$ws.Range("P6").Value = 8.594181508245255
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = 3
$ws.Range("U6").Value = 45
$ws.Range("V6").Value = 2.512658227848101
$ws.Range("W6").Value = 4.48098251948876
$ws.Range("X6").Value = 0
$ws.Range("Y6").Value = 0
$ws.Range("Z6").Value = 0
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 17
$ws.Range("AC6").Value = 1.069620253164557
$ws.Range("AD6").Value = 3.073939798870591
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 0
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 2
$ws.Range("AJ6").Value = 0.120253164556962
$ws.Range("AK6").Value = 0.3447123472852491
$ws.Range("AL6").Value = 0.6666666666666666
$ws.Range("AM6").Value = 1.5
$ws.Range("AN6").Value = 2
$ws.Range("AO6").Value = 3
$ws.Range("AP6").Value = 6
$ws.Range("AQ6").Value = 2.290964753656293
$ws.Range("AR6").Value = 0.9224471218358076
$ws.Range("AS6").Value = 9.333333333333334
$ws.Range("AT6").Value = 13.57142857142857
$ws.Range("AU6").Value = 15.83333333333333
$ws.Range("AV6").Value = 19
$ws.Range("AW6").Value = 34
$ws.Range("AX6").Value = 16.53393065646352
$ws.Range("AY6").Value = 4.521254120415229

# Row 11 (new last row, beyond the old A1:AY9 used range, which auto-extends
# the sheet's dimension to A1:AY11): P naive3b
$ws.Range("A11").Value = "P naive3b"
$ws.Range("B11").Value = 254
$ws.Range("C11").Value = 49
$ws.Range("D11").Value = 95
$ws.Range("E11").Value = 95
$ws.Range("F11").Value = 97
$ws.Range("G11").Value = 428
$ws.Range("H11").Value = 152.3503937007874
$ws.Range("I11").Value = 129.265450066412
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 4
$ws.Range("L11").Value = 6
$ws.Range("M11").Value = 7
$ws.Range("N11").Value = 35
$ws.Range("O11").Value = 7.377952755905512
$ws.Range("P11").Value = 5.779107982453519
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = 1
$ws.Range("S11").Value = 2
$ws.Range("T11").Value = 5
$ws.Range("U11").Value = 17
$ws.Range("V11").Value = 3.216535433070866
$ws.Range("W11").Value = 3.275909523513953
$ws.Range("X11").Value = 0
$ws.Range("Y11").Value = 0
$ws.Range("Z11").Value = 0
$ws.Range("AA11").Value = 0
$ws.Range("AB11").Value = 19
$ws.Range("AC11").Value = 1.803149606299213
$ws.Range("AD11").Value = 3.905078015807789
$ws.Range("AE11").Value = 0
$ws.Range("AF11").Value = 0
$ws.Range("AG11").Value = 0
$ws.Range("AH11").Value = 0
$ws.Range("AI11").Value = 2
$ws.Range("AJ11").Value = 0.2007874015748032
$ws.Range("AK11").Value = 0.4299087251214076
$ws.Range("AL11").Value = 0
$ws.Range("AM11").Value = 1.333333333333333
$ws.Range("AN11").Value = 2
$ws.Range("AO11").Value = 3
$ws.Range("AP11").Value = 4.5
$ws.Range("AQ11").Value = 2.192074674876167
$ws.Range("AR11").Value = 1.038702278113778
$ws.Range("AS11").Value = 9.5
$ws.Range("AT11").Value = 13.57142857142857
$ws.Range("AU11").Value = 15.83333333333333
$ws.Range("AV11").Value = 19.4
$ws.Range("AW11").Value = 94
$ws.Range("AX11").Value = 18.73246766019657
$ws.Range("AY11").Value = 10.16469840167657
